# Update timestamped test email addresses (appear on "UsuariosRegistro"
# and some are reused on "LoginData") to reflect the new timestamp.
# Old timestamp: 20251109_022039 -> New timestamp: 20251109_024842

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_022039"
$newStamp = "20251109_024842"

function Update-StampedEmails($range) {
    foreach ($cell in $range.Cells) {
        $txt = $cell.Text
        if ($txt -like "*$oldStamp*") {
            $cell.Value = $txt.Replace($oldStamp, $newStamp)
        }
    }
}

# Main user registration data: emails live in column C, rows 2-6
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
Update-StampedEmails $wsUsuarios.Range("C2:C6")

# Login test data reuses some of the same e-mail addresses in column A
$wsLogin = $wb.Worksheets.Item("LoginData")
Update-StampedEmails $wsLogin.Range("A2:A6")
